$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 40497
$ws.Range("A3").Value = 40678
$ws.Range("A4").Value = 40862
$ws.Range("A5").Value = 41044
$ws.Range("A6").Value = 41228
$ws.Range("A7").Value = 41409
$ws.Range("A8").Value = 41593
$ws.Range("A9").Value = 41774
$ws.Range("A10").Value = 41958
$ws.Range("A11").Value = 42139
$ws.Range("A12").Value = 42323
$ws.Range("A13").Value = 42505
$ws.Range("A14").Value = 42689
$ws.Range("A15").Value = 42870
$ws.Range("A16").Value = 43054
$ws.Range("A17").Value = 43146
$ws.Range("A18").Value = 43235
$ws.Range("A19").Value = 43327
$ws.Range("A20").Value = 43419
$ws.Range("A21").Value = 43511
$ws.Range("A22").Value = 43600
$ws.Range("A23").Value = 43692
$ws.Range("A24").Value = 43784
$ws.Range("A25").Value = 43876
$ws.Range("A26").Value = 43966
$ws.Range("A27").Value = 44058
$ws.Range("A28").Value = 44150
$ws.Range("A29").Value = 44242
$ws.Range("A30").Value = 44331
$ws.Range("A31").Value = 44423
$ws.Range("A32").Value = 44515
$ws.Range("A33").Value = 44607
$ws.Range("A34").Value = 44696
$ws.Range("A35").Value = 44788
$ws.Range("A36").Value = 44880
$ws.Range("A37").Value = 44972
$ws.Range("A38").Value = 45061
$ws.Range("A39").Value = 45153
$ws.Range("A40").Value = 45245
$ws.Range("A41").Value = 45337
$ws.Range("A42").Value = 45427
$ws.Range("A43").Value = 45519
$ws.Range("A44").Value = 45611
$ws.Range("A45").Value = 45703
$ws.Range("A46").Value = 45792
$ws.Range("A47").Value = 45884
